# Add extra columns (WIN, TOP4, TOP5, RELEGATION) ahead of ExpPoints,
# shift ExpPoints from column C to column G, update the Team order / values
# for the new matchday-8 prediction, preparing for the Monte Carlo simulation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Propagate the header styling (bold font, boxed border, centred text)
# from the existing C1 ("ExpPoints") cell onto the new D1:G1 cells by
# copying the formatted cell across, then overwrite the copied text.
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("C1").Copy($ws.Range("E1"))
$ws.Range("C1").Copy($ws.Range("F1"))
$ws.Range("C1").Copy($ws.Range("G1"))

$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "RELEGATION"
$ws.Range("G1").Value = "ExpPoints"

# --- Updated team order / expected points ------------------------------
$teams = @(
    "Arsenal",
    "Manchester City",
    "Liverpool",
    "Chelsea",
    "Crystal Palace",
    "Aston Villa",
    "AFC Bournemouth",
    "Brighton & Hove Albion",
    "Newcastle United",
    "Manchester United",
    "Tottenham Hotspur",
    "Brentford",
    "Everton",
    "Fulham",
    "Sunderland",
    "Nottingham Forest",
    "Leeds United",
    "West Ham United",
    "Burnley",
    "Wolverhampton Wanderers"
)

$expPoints = @(
    80.72983931743445,
    73.01948554978178,
    70.01200598378286,
    62.82593626625423,
    58.62605350972558,
    58.14948380623721,
    55.81627037770095,
    55.60172351896105,
    54.43896264218951,
    52.92965991313606,
    52.30056981730057,
    50.61051401169723,
    45.36584335703066,
    42.91340272334515,
    41.66784422819364,
    37.23513749964833,
    35.91542427582408,
    34.52745126016377,
    34.13497402429287,
    29.82945707020081
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 2).Value = $teams[$i]

    # Clear out the old ExpPoints value, columns C-F become blank
    # placeholder cells for the upcoming Monte Carlo percentages.
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 6).Value = ""

    $ws.Cells.Item($row, 7).Value = $expPoints[$i]
}
